$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data set (header + 30 data rows). Rows 2-4 are newly
# inserted at the top of the data (pushing the former rows 2-21 down by
# three), and rows 25-31 are newly appended at the bottom.
$allRows = @(
    @(-0.03444957733154287, 0.01032954454422018, 0.2995094756285349),
    @(-0.07984677950541165, -0.1310317392150563, 0.2751492535074553),
    @(-0.09086012840270966, -0.02373796701431304, 0.2084361910820008),
    @(0.005503296852111615, -0.1945102314154306, 0.2710902802646161),
    @(-0.01686970392863169, -0.689376473426814, -0.4833486676216042),
    @(0.7827135920524584, -2.646403878927229, -4.282532811164854),
    @(-1.04310808579127, 0.482522944609324, -1.852447708447774),
    @(-2.88519119222959, -0.1530523101488742, 0.5539670785268116),
    @(-2.359475612640381, -0.4427683353424072, 1.851428031921387),
    @(-3.12837818264961, 1.37843454877535, 1.588190158208212),
    @(-9.990091621875759, 0.5772410134474441, 7.240908702214556),
    @(3.089871168136629, 1.149596989154817, -7.327601671218906),
    @(-1.657969951629632, -0.6963444749514238, -3.005745251973476),
    @(-2.83851852019628, -0.4919258654117603, -0.4817840854326922),
    @(-1.064056754112247, 0.9830632209777808, -6.221896290779103),
    @(-0.642586628595987, -0.08318759004274168, -0.6120627919833304),
    @(-2.159881273905432, -1.468431472778319, -0.03336954116820623),
    @(-3.000638425350189, -0.02016200125217193, -1.370798766613007),
    @(-4.049820343653357, 1.876313845316566, 1.094406366348257),
    @(-7.227878630161286, 1.69673324127992, 3.92280129591624),
    @(-2.149182796478235, -1.535205185413374, 0.8277400732040205),
    @(3.039104461669923, -1.047793030738832, -4.546324898799259),
    @(-1.214001258214315, 1.53115118543307, -4.196547468503316),
    @(-2.15471959114075, 0.2546487897634502, -0.9123589172959333),
    @(0.6330445607503264, 0.5185447335243235, 0.3759170770645095),
    @(-0.2923502524693779, -0.1014364187916107, 1.366777941584587),
    @(-1.484014511108395, -0.6635967791080455, 1.253775984048842),
    @(0.2468122641245533, 0.4217223922411623, 0.5185143599907569),
    @(-0.3225175539652509, -0.6644023060798641, -0.07338536779085762),
    @(-0.4933383464813202, -0.1615586355328552, 0.263843480497599)
)

$r = 2
foreach ($row in $allRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}
